# Fig18.xlsx - 2017-02-13 snapshot (STEO January 2017 -> February 2017 roll)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig18")

# --- Title / source caption text (shared strings A2 and B30) ---
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("B30").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Natural gas production table (rows 27-29), columns D:H are raw inputs,
#     columns J:M hold =E-D shared formulas that recompute automatically. ---
$ws.Range("F27").Value = 3.348185

$ws.Range("F28").Value = 74.033877667070001
$ws.Range("G28").Value = 75.788586886190004
$ws.Range("H28").Value = 80.398135766030009

$ws.Range("F29").Value = 1.8171688257
$ws.Range("G29").Value = 0.51082380630000002
$ws.Range("H29").Value = -1.0501675402999999

# --- History / forecast series used by the chart (rows 35-82) ---
# History (column C) extends one more month; forecast (column D) shifts out.
$ws.Range("C55").Value = 76.756022133000002
$ws.Range("C56").Value = 75.832430290000005
$ws.Range("C57").Value = 76.995246366999993

$ws.Range("C58").Value = 76.778400000000005
$ws.Range("D58").Value = "#N/A"

$ws.Range("C59").Value = 76.154250000000005
$ws.Range("D59").Value = 76.154250000000005

$ws.Range("D60").Value = 76.884979999999999
$ws.Range("D61").Value = 77.70317
$ws.Range("D62").Value = 77.872829999999993
$ws.Range("D63").Value = 78.358279999999993
$ws.Range("D64").Value = 78.595830000000007
$ws.Range("D65").Value = 79.320250000000001
$ws.Range("D66").Value = 79.982290000000006
$ws.Range("D67").Value = 80.357089999999999
$ws.Range("D68").Value = 80.819050000000004
$ws.Range("D69").Value = 81.173590000000004
$ws.Range("D70").Value = 81.425460000000001
$ws.Range("D71").Value = 81.957369999999997
$ws.Range("D72").Value = 82.583839999999995
$ws.Range("D73").Value = 83.086860000000001
$ws.Range("D74").Value = 83.305459999999997
$ws.Range("D75").Value = 83.476060000000004
$ws.Range("D76").Value = 83.55641
$ws.Range("D77").Value = 83.745369999999994
$ws.Range("D78").Value = 84.062889999999996
$ws.Range("D79").Value = 84.142120000000006
$ws.Range("D80").Value = 84.424880000000002
$ws.Range("D81").Value = 84.725930000000005
$ws.Range("D82").Value = 84.977080000000001
